# Auto-generated Excel COM-interop edit script
# Applies cell-value updates to the Unicorn_Profits workbook sheets
# (data refresh from the scheduled market-data runner).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 70
$ws.Range("H70").Value = 3525.6296
$ws.Range("I70").Value = 8349
$ws.Range("J70").Value = 1494.7368
$ws.Range("K70").Value = 25047
$ws.Range("L70").Value = 4484.2104
$ws.Range("M70").Value = -24777
$ws.Range("N70").Value = -5024.2104

# Row 73
$ws.Range("H73").Value = 3525.6296
$ws.Range("I73").Value = 8349
$ws.Range("J73").Value = 1494.7368
$ws.Range("K73").Value = 25047
$ws.Range("L73").Value = 4484.2104
$ws.Range("M73").Value = -24111
$ws.Range("N73").Value = -6356.2104

# Row 82
$ws.Range("H82").Value = 447.33334
$ws.Range("I82").Value = 447.33334
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1342.00002
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -936.0000199999999
$ws.Range("N82").ClearContents()

# Row 85
$ws.Range("H85").Value = 447.33334
$ws.Range("I85").Value = 447.33334
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1342.00002
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 61.99998000000005
$ws.Range("N85").ClearContents()

# Row 132
$ws.Range("H132").Value = 2760.695
$ws.Range("I132").Value = 1242.0968
$ws.Range("J132").Value = 7468.35
$ws.Range("K132").Value = 3726.2904
$ws.Range("L132").Value = 22405.05
$ws.Range("M132").Value = -1196.2904
$ws.Range("N132").Value = -27465.05


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 1589.8
$ws.Range("I32").Value = 1383.299
$ws.Range("J32").Value = 8266.666999999999
$ws.Range("K32").Value = 1383.299
$ws.Range("L32").Value = 8266.666999999999
$ws.Range("M32").Value = -1096.299
$ws.Range("N32").Value = -8840.666999999999

# Row 113
$ws.Range("H113").Value = 31010.5
$ws.Range("J113").Value = 31010.5
$ws.Range("L113").Value = 31010.5
$ws.Range("N113").Value = -39688.5

# Row 122
$ws.Range("H122").Value = 4354.1816
$ws.Range("I122").Value = 3905.111
$ws.Range("J122").Value = 6375
$ws.Range("K122").Value = 11715.333
$ws.Range("L122").Value = 19125
$ws.Range("M122").Value = -9265.332999999999
$ws.Range("N122").Value = -24025


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 54
$ws.Range("H54").Value = 5986.0625
$ws.Range("I54").Value = 2141.7
$ws.Range("J54").Value = 12393.333
$ws.Range("K54").Value = 2141.7
$ws.Range("L54").Value = 12393.333
$ws.Range("M54").Value = -1657.7
$ws.Range("N54").Value = -13361.333


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 3003.0908
$ws.Range("I31").Value = 1710.4517
$ws.Range("J31").Value = 6085.5386
$ws.Range("K31").Value = 1710.4517
$ws.Range("L31").Value = 6085.5386
$ws.Range("M31").Value = -1415.4517
$ws.Range("N31").Value = -6675.5386

# Row 34
$ws.Range("H34").Value = 3003.0908
$ws.Range("I34").Value = 1710.4517
$ws.Range("J34").Value = 6085.5386
$ws.Range("K34").Value = 1710.4517
$ws.Range("L34").Value = 6085.5386
$ws.Range("M34").Value = -1508.4517
$ws.Range("N34").Value = -6489.5386

# Row 36
$ws.Range("H36").Value = 4800
$ws.Range("I36").Value = 4800
$ws.Range("K36").Value = 4800
$ws.Range("M36").Value = -4412

# Row 40
$ws.Range("H40").Value = 4800
$ws.Range("I40").Value = 4800
$ws.Range("K40").Value = 4800
$ws.Range("M40").Value = -4640

# Row 58
$ws.Range("H58").Value = 2344.9473
$ws.Range("I58").Value = 2556.25
$ws.Range("J58").Value = 1887.125
$ws.Range("K58").Value = 2556.25
$ws.Range("L58").Value = 1887.125
$ws.Range("M58").Value = -2353.25
$ws.Range("N58").Value = -2293.125

# Row 136
$ws.Range("H136").Value = 2344.9473
$ws.Range("I136").Value = 2556.25
$ws.Range("J136").Value = 1887.125
$ws.Range("K136").Value = 7668.75
$ws.Range("L136").Value = 5661.375
$ws.Range("M136").Value = -5118.75
$ws.Range("N136").Value = -10761.375


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 4
$ws.Range("H4").Value = 1011.82355
$ws.Range("I4").Value = 112.375
$ws.Range("J4").Value = 1811.3334
$ws.Range("K4").Value = 337.125
$ws.Range("L4").Value = 5434.0002
$ws.Range("M4").Value = -225.125
$ws.Range("N4").Value = -5658.0002

# Row 5
$ws.Range("H5").Value = 571.4783
$ws.Range("I5").Value = 531.2632
$ws.Range("J5").Value = 762.5
$ws.Range("K5").Value = 1593.7896
$ws.Range("L5").Value = 2287.5
$ws.Range("M5").Value = -1481.7896
$ws.Range("N5").Value = -2511.5

# Row 131
$ws.Range("H131").Value = 1143.9294
$ws.Range("I131").Value = 1449.9
$ws.Range("J131").Value = 1103.1333
$ws.Range("K131").Value = 4349.700000000001
$ws.Range("L131").Value = 3309.3999
$ws.Range("M131").Value = 690.2999999999993
$ws.Range("N131").Value = -13389.3999

# Row 132
$ws.Range("H132").Value = 5846.0586
$ws.Range("I132").Value = 3387.111
$ws.Range("J132").Value = 8612.375
$ws.Range("K132").Value = 30483.999
$ws.Range("L132").Value = 77511.375
$ws.Range("M132").Value = -27953.999
$ws.Range("N132").Value = -82571.375

# Row 135
$ws.Range("H135").Value = 571.4783
$ws.Range("I135").Value = 531.2632
$ws.Range("J135").Value = 762.5
$ws.Range("K135").Value = 4781.3688
$ws.Range("L135").Value = 6862.5
$ws.Range("M135").Value = -2246.3688
$ws.Range("N135").Value = -11932.5


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 122
$ws.Range("H122").Value = 1210.7778
$ws.Range("I122").Value = 809.4
$ws.Range("J122").Value = 1712.5
$ws.Range("K122").Value = 2428.2
$ws.Range("L122").Value = 5137.5
$ws.Range("M122").Value = 21.80000000000018
$ws.Range("N122").Value = -10037.5


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 2
$ws.Range("H2").Value = 250002620
$ws.Range("J2").Value = 250002620
$ws.Range("L2").Value = 250002620
$ws.Range("N2").Value = -250002844

# Row 103
$ws.Range("H103").Value = 38000
$ws.Range("J103").Value = 38000
$ws.Range("L103").Value = 38000
$ws.Range("N103").Value = -40344

# Row 106
$ws.Range("H106").Value = 19049
$ws.Range("J106").Value = 19049
$ws.Range("L106").Value = 19049
$ws.Range("N106").Value = -21573


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 2
$ws.Range("H2").Value = 28644.666
$ws.Range("I2").Value = 9954
$ws.Range("J2").Value = 47335.332
$ws.Range("K2").Value = 9954
$ws.Range("L2").Value = 47335.332
$ws.Range("M2").Value = -9842
$ws.Range("N2").Value = -47559.332

# Row 103
$ws.Range("H103").Value = 26650.5
$ws.Range("J103").Value = 26650.5
$ws.Range("L103").Value = 26650.5
$ws.Range("N103").Value = -28994.5

# Row 132
$ws.Range("H132").Value = 1733.0361
$ws.Range("I132").Value = 1118.44
$ws.Range("J132").Value = 2664.2424
$ws.Range("K132").Value = 3355.32
$ws.Range("L132").Value = 7992.7272
$ws.Range("M132").Value = -825.3200000000002
$ws.Range("N132").Value = -13052.7272


Write-Output "Applied Unicorn_Profits data refresh edits."